# Update the Queens (NY)_B team-specific transition probability matrix.
# Rows/cols correspond to the states listed in column A / row 1 (Af0..Af3, Ai0..Ai3, Ar0, Bf0..Bf3, Bi0..Bi3, Br0).
# Values below are the simulated transition probabilities produced by the (sped-up) game simulator.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2307692307692308
$ws.Range("C2").Value = 0.4615384615384616
$ws.Range("S2").Value = 0.3076923076923077

# Row 3
$ws.Range("P3").Value = 0.5
$ws.Range("S3").Value = 0.5

# Row 4
$ws.Range("P4").Value = 1

# Row 5
$ws.Range("P5").Value = 1

# Row 6
$ws.Range("B6").Value = 0.1428571428571428
$ws.Range("J6").Value = 0.2857142857142857
$ws.Range("Q6").Value = 0.1428571428571428
$ws.Range("S6").Value = 0.4285714285714285

# Row 7
$ws.Range("B7").Value = 0.1428571428571428
$ws.Range("F7").Value = 0.1428571428571428
$ws.Range("J7").Value = 0.1428571428571428
$ws.Range("Q7").Value = 0.1428571428571428
$ws.Range("S7").Value = 0.4285714285714285

# Row 8
$ws.Range("B8").Value = 0.1578947368421053
$ws.Range("F8").Value = 0.05263157894736842
$ws.Range("J8").Value = 0.05263157894736842
$ws.Range("Q8").Value = 0.05263157894736842
$ws.Range("R8").Value = 0.05263157894736842
$ws.Range("S8").Value = 0.631578947368421

# Row 9
$ws.Range("R9").Value = 0.3333333333333333
$ws.Range("S9").Value = 0.6666666666666666

# Row 10
$ws.Range("B10").Value = 0.1176470588235294
$ws.Range("D10").Value = 0.0392156862745098
$ws.Range("E10").Value = 0.0196078431372549
$ws.Range("F10").Value = 0.0392156862745098
$ws.Range("J10").Value = 0.07843137254901961
$ws.Range("O10").Value = 0.0196078431372549
$ws.Range("Q10").Value = 0.1568627450980392
$ws.Range("R10").Value = 0.09803921568627451
$ws.Range("S10").Value = 0.4313725490196079

# Row 11
$ws.Range("G11").Value = 0.07142857142857142
$ws.Range("J11").Value = 0.07142857142857142
$ws.Range("K11").Value = 0.07142857142857142
$ws.Range("L11").Value = 0.7857142857142857

# Row 12
$ws.Range("G12").Value = 0.5833333333333334
$ws.Range("J12").Value = 0.3333333333333333
$ws.Range("L12").Value = 0.08333333333333333

# Row 13
$ws.Range("J13").Value = 1

# Row 15
$ws.Range("H15").Value = 0.125
$ws.Range("J15").Value = 0.625
$ws.Range("K15").Value = 0.125
$ws.Range("S15").Value = 0.125

# Row 16
$ws.Range("H16").Value = 0.2
$ws.Range("J16").Value = 0.4
$ws.Range("K16").Value = 0.2
$ws.Range("M16").Value = 0.2

# Row 17
$ws.Range("F17").Value = 0.09090909090909091
$ws.Range("H17").Value = 0.1818181818181818
$ws.Range("J17").Value = 0.6363636363636364
$ws.Range("S17").Value = 0.09090909090909091

# Row 18
$ws.Range("I18").Value = 0.1428571428571428
$ws.Range("J18").Value = 0.5714285714285714
$ws.Range("K18").Value = 0.1428571428571428
$ws.Range("S18").Value = 0.1428571428571428

# Row 19
$ws.Range("F19").Value = 0.01754385964912281
$ws.Range("H19").Value = 0.2631578947368421
$ws.Range("I19").Value = 0.03508771929824561
$ws.Range("J19").Value = 0.3157894736842105
$ws.Range("K19").Value = 0.1578947368421053
$ws.Range("M19").Value = 0.01754385964912281
$ws.Range("O19").Value = 0.1052631578947368
$ws.Range("S19").Value = 0.08771929824561403
